# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> used by the Slide Master (was "Integral" / "Red Violet")
#   ppt/theme/theme2.xml -> used by the Notes Master  (was "Office Theme" / "Office")
#
# The authored edit swaps the two colour palettes: the Slide Master's theme
# becomes the stock "Office" palette, and the Notes Master's theme becomes the
# "Red Violet" palette that used to live on the Slide Master. Font/format
# schemes are identical between the two themes and are left untouched.
#
# Theme colours are edited through ThemeColorScheme.Colors(index).RGB, the
# PowerPoint COM surface for per-slot theme colour edits. The twelve slots,
# in COM index order, map onto the OOXML <a:clrScheme> children as:
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink   12 folHlink

function Set-ThemeColors($themeColorScheme, $hexColors) {
    for ($i = 0; $i -lt $hexColors.Count; $i++) {
        $hex = $hexColors[$i]
        $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
        $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
        $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
        $oleRgb = $r + ($g * 256) + ($b * 65536)
        $themeColorScheme.Colors($i + 1).RGB = $oleRgb
    }
}

$p = $ppt.ActivePresentation

# Order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")
$redVioletColors = @("000000", "FFFFFF", "454551", "D8D9DC", "E32D91", "C830CC", "4EA6DC", "4775E7", "8971E1", "D54773", "6B9F25", "8C8C8C")

# Slide Master's theme (theme1.xml): Integral/Red Violet -> Office Theme/Office
$slideMasterTheme = $p.SlideMaster.Theme.ThemeColorScheme
Set-ThemeColors $slideMasterTheme $officeColors

# Notes Master's theme (theme2.xml): Office Theme/Office -> Integral/Red Violet
$notesMasterTheme = $p.NotesMaster.Theme.ThemeColorScheme
Set-ThemeColors $notesMasterTheme $redVioletColors
